$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-02-09 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-10 Monday", 2) | Out-Null

# Update the multiplication answers in the table, cell by cell (row, col)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "47×32=1504"
$t.Cell(1, 2).Range.Text = "56×77=4312"
$t.Cell(1, 3).Range.Text = "29×56=1624"
$t.Cell(1, 4).Range.Text = "85×44=3740"
$t.Cell(1, 5).Range.Text = "11×43=473"

$t.Cell(5, 1).Range.Text = "60×88=5280"
$t.Cell(5, 2).Range.Text = "71×58=4118"
$t.Cell(5, 3).Range.Text = "23×91=2093"
$t.Cell(5, 4).Range.Text = "55×24=1320"
$t.Cell(5, 5).Range.Text = "91×75=6825"

$t.Cell(10, 1).Range.Text = "16×69=1104"
$t.Cell(10, 2).Range.Text = "31×24=744"
$t.Cell(10, 3).Range.Text = "71×81=5751"
$t.Cell(10, 4).Range.Text = "15×85=1275"
$t.Cell(10, 5).Range.Text = "56×66=3696"

$t.Cell(15, 1).Range.Text = "44×98=4312"
$t.Cell(15, 2).Range.Text = "27×76=2052"
$t.Cell(15, 3).Range.Text = "84×69=5796"
$t.Cell(15, 4).Range.Text = "72×31=2232"
$t.Cell(15, 5).Range.Text = "93×84=7812"

$t.Cell(20, 1).Range.Text = "42×24=1008"
$t.Cell(20, 2).Range.Text = "12×78=936"
$t.Cell(20, 3).Range.Text = "94×13=1222"
$t.Cell(20, 4).Range.Text = "25×31=775"
$t.Cell(20, 5).Range.Text = "63×19=1197"
